$d = $word.ActiveDocument

# Locate the two "Programa" (not "Programa resumido") body paragraphs that
# contain the run-on text without spaces between sentences - these are the
# paragraphs the diff targets (the ones right after the "Programa" heading).
$ptParaIndex = -1
$enParaIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Introdu*nanotecnologia.Spintr*") {
        $ptParaIndex = $i
    }
    if ($t -like "Introduction to nanotechnology.Metal*") {
        $enParaIndex = $i
    }
    $i = $i + 1
}

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$ptBody = '<w:p><w:r><w:t>Introdução à nanotecnologia.</w:t><w:br/><w:t xml:space="preserve">Spintrônica de metais. Spintrônica de semicondutores </w:t><w:br/><w:t>Dispositivos da spintrônica.</w:t><w:br/><w:t>Introdução à computação clássica. Introdução à computação quântica. Algoritmos quânticos.</w:t><w:br/><w:t xml:space="preserve">Decoerência. Pontos quânticos. </w:t><w:br/><w:t>Transistor de Kane.</w:t><w:br/><w:t>Introdução a modelos da consciência: o cérebro é um computador quântico?</w:t></w:r></w:p>'

$enBody = '<w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Introduction to nanotechnology.</w:t><w:br/><w:t>Metal spintronics. Semiconductor Spintronics</w:t><w:br/><w:t>Spintronics devices.</w:t><w:br/><w:t>Introduction to classical computing. Introduction to quantum computing. Quantum Algorithms.</w:t><w:br/><w:t>decoherence. Quantum Dots.</w:t><w:br/><w:t>Kane transistor.</w:t><w:br/><w:t>Introduction to models of consciousness: is the brain a quantum computer?</w:t></w:r></w:p>'

$ptRange = $d.Paragraphs($ptParaIndex).Range
$ptRange.InsertXML($pkgHeader + $ptBody + $pkgFooter)

$enRange = $d.Paragraphs($enParaIndex).Range
$enRange.InsertXML($pkgHeader + $enBody + $pkgFooter)
